$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = "2026-01-14 18:39:14"
}
